$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value()

    if ($null -ne $val -and $val -ne "") {
        $parts = $val -split ', '
        $hasSystem = $false
        $others = @()

        foreach ($p in $parts) {
            if ($p.Equals("System")) {
                $hasSystem = $true
            } else {
                $others += $p
            }
        }

        if ($hasSystem -and $others.Count -gt 0) {
            $newParts = @("System") + $others
            $newVal = [string]::Join(", ", $newParts)
            $cell.Value = $newVal
        }
    }
}
